# Added Liga MX to daily picks
# Appends 18 Liga MX clubs (Name / FBRef-short-name / Bovada-name) to the
# bottom of the Name/FBRef/Bovada lookup table on Sheet1, starting at row 304.
#
# Columns: A = "Name" (long/site display form), B = "FBRef" (short form),
#          C = "Bovada" (same text as column A for every row on this sheet).
#
# The rows were originally authored by pasting/typing column B first (sorted
# alphabetically, top to bottom) and then column A/C afterwards (also
# alphabetically, but landing on different rows) -- we reproduce the same
# two-pass fill order so new entries line up the way they were written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, Name (col A / C), FBRef short name (col B)
$ligaMx = @(
    [PSCustomObject]@{Row=304; Name="Club America";        Short="América"},
    [PSCustomObject]@{Row=305; Name="Atlas FC";            Short="Atlas"},
    [PSCustomObject]@{Row=306; Name="Atletico San Luis";   Short="Atlético"},
    [PSCustomObject]@{Row=307; Name="Cruz Azul";           Short="Cruz Azul"},
    [PSCustomObject]@{Row=308; Name="FC Juarez";           Short="FC Juárez"},
    [PSCustomObject]@{Row=309; Name="Chivas Guadalajara";  Short="Guadalajara"},
    [PSCustomObject]@{Row=310; Name="Club Leon";           Short="León"},
    [PSCustomObject]@{Row=311; Name="Mazatlan FC";         Short="Mazatlán"},
    [PSCustomObject]@{Row=312; Name="CF Monterrey";        Short="Monterrey"},
    [PSCustomObject]@{Row=313; Name="Necaxa";              Short="Necaxa"},
    [PSCustomObject]@{Row=314; Name="CF Pachuca";          Short="Pachuca"},
    [PSCustomObject]@{Row=315; Name="Puebla FC";           Short="Puebla"},
    [PSCustomObject]@{Row=316; Name="Queretaro FC";        Short="Querétaro"},
    [PSCustomObject]@{Row=317; Name="Santos Laguna";       Short="Santos"},
    [PSCustomObject]@{Row=318; Name="Tijuana";             Short="Tijuana"},
    [PSCustomObject]@{Row=319; Name="Deportivo Toluca FC"; Short="Toluca"},
    [PSCustomObject]@{Row=320; Name="Tigres UANL";         Short="UANL"},
    [PSCustomObject]@{Row=321; Name="Pumas UNAM";          Short="UNAM"}
)

# Pass 1: fill column B (FBRef) top-to-bottom -- this list is already in
# alphabetical order by the short name.
foreach ($entry in $ligaMx) {
    $ws.Cells.Item($entry.Row, 2).Value = $entry.Short
}

# Pass 2: fill column A (Name), visiting the rows in alphabetical order of
# the long name (the order those values were originally typed/pasted in).
$byName = $ligaMx | Sort-Object -Property Name
foreach ($entry in $byName) {
    $ws.Cells.Item($entry.Row, 1).Value = $entry.Name
}

# Pass 3: fill column C (Bovada) -- identical text to column A, row by row.
foreach ($entry in $ligaMx) {
    $ws.Cells.Item($entry.Row, 3).Value = $entry.Name
}

# Leave the selection where the author ended up after typing the last row.
$ws.Range("E304").Select() | Out-Null
